$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 95941.66286295939
$ws.Range("B3").Value = 95714.4546870803
$ws.Range("B4").Value = 95796.14810636583
$ws.Range("B5").Value = 95651.19403706689
$ws.Range("B6").Value = 95952.93469409621
$ws.Range("B7").Value = 95899.98587160894
$ws.Range("B8").Value = 95787.18861133493
$ws.Range("B9").Value = 96010.12640343809
$ws.Range("B10").Value = 95702.49673202615
$ws.Range("B11").Value = 95770.74633360193
$ws.Range("B12").Value = 95625.41423583131
$ws.Range("B13").Value = 95984.49127048081
$ws.Range("B14").Value = 95809.35272688311
$ws.Range("B15").Value = 95956.61864088102
$ws.Range("B16").Value = 95745.1099113618
$ws.Range("B17").Value = 95788.00999194197
$ws.Range("B18").Value = 96213.36815600557
$ws.Range("B19").Value = 95971.94085414987
$ws.Range("B20").Value = 95671.19817351598
$ws.Range("B21").Value = 95656.88068761754
$ws.Range("B22").Value = 95641.21638272796
$ws.Range("B23").Value = 95871.59656191245
$ws.Range("B24").Value = 96001.99221058287
$ws.Range("B25").Value = 96069.84850926674
$ws.Range("B26").Value = 95767.5107146684
$ws.Range("B27").Value = 95928.56803652969
$ws.Range("B28").Value = 95913.82535589579
$ws.Range("B29").Value = 95742.96099919418
$ws.Range("B30").Value = 95684.94771241832
$ws.Range("B31").Value = 95870.20703733548
$ws.Range("B32").Value = 95953.78291700241
$ws.Range("B33").Value = 95803.18162771959
$ws.Range("B34").Value = 95786.86558448516
$ws.Range("B35").Value = 95936.17539618588
$ws.Range("B36").Value = 95892.96180499597
$ws.Range("B37").Value = 95800.7238248724
$ws.Range("B38").Value = 95987.95446265939
$ws.Range("B39").Value = 95750.26768734891
